$d = $word.ActiveDocument

# --- Paragraph 6: Objetivos summary (PT) ---
$d.Paragraphs.Item(6).Range.Text = "Cinética Química, Termodinâmica e Equilíbrio, Equilíbrio Químico, Eletroquímica."

# --- Paragraph 7: Objetivos summary (EN, italic) ---
$d.Paragraphs.Item(7).Range.Text = "Chemical Kinetics, Thermodynamics and Equilibrium, Chemical Equilibrium, Eletrochemistry."

# --- Paragraph 9: Docentes list bullet -> long PT objectives + Programa detail text ---
$d.Paragraphs.Item(9).Range.Text = "Fazer previsões e explicar acerca dos efeitos sobre as velocidades das reações em vista de: catálise, variação da temperatura, geometria de colisão e concentração dos reagentes. Sugerir leis de velocidade de reação de posse de dados apropriados acerca dos efeitos de concentração, estudar as reações de ordens zero, um e dois. Estudar a aplicação da equação de Arrhenius. Ampliar o entendimento do sentido espontâneo das reações químicas. Entender a função termodinâmica entropia e sua relação com as três leis da termodinâmica. Entender o sentido de energia livre como uma referência para o grau de afastamento entre o sistema e seu estado de equilíbrio. Utilizar a variação da energia livre padrão como uma ferramenta para calcular a constante de equilíbrio para determinado processo. Examinar o conceito de equilíbrio e definir a constante de equilíbrio. Aprender a escrever as expressões das constantes de equilíbrio para reações homogêneas e heterogêneas e interpretar o sentido do quociente reacional. Dar a expressão do produto de solubilidade para um sal e calcular os produtos de solubilidade a partir de solubilidades determinadas experimentalmente e prever se deve ou não ocorrer precipitação. Aplicar os critérios de precipitação. Determinar os efeitos do íon comum. Calcular o pH de uma solução aquosa para sistemas envolvendo ácido ou base fortes ou pelo menos um ácido fraco ou uma base fraca. Entender o funcionamento de soluções-tampão. Estudar as reações envolvendo a formação de íons complexos a partir de espécies muito solúveis e muito pouco solúveis. Realizar uma representação simples para um sistema químico de uma pilha eletrolítica ou “galvânica” observando as convenções-padrão para identificar os eletrodos. Aplicar os princípios de estequiometria aos processos eletroquímicos usando equações balanceadas de semi-reações e o valor da constante de Faraday. Calcular potenciais-padrão de pilhas usando valores tabelados. Energia livre de Gibbs e a equação de Nernst." + [char]11 + "Cinética Química: Energia de ativação e catalisadores. Efeito da concentração dos reagentes e temperatura sobre a velocidade de reações químicas. Equações diferenciais de velocidade de reação. Leis de velocidade para reações de 1ª e 2ª ordens. Equação de Arrhenius. " + [char]11 + "Termodinâmica e equilíbrio: Primeira Lei da Termodinâmica. Calor de reação e energia interna. Entalpia e variação de Entalpia. Entropia. Segunda Lei da Termodinâmica. Entropia e Desordem. Terceira Lei da Termodinâmica. Variação de Entropia numa reação. Energia Livre Padrões e Critério de Espontaneidade. Relação entre variação de energia livre padrão e a constante de Equilíbrio. " + [char]11 + "Equilíbrio Químico: Natureza do equilíbrio químico. Quociente de reação e constante de equilíbrio. Efeito da concentração de reagentes e temperatura sobre o equilíbrio (princípio de Le Chatelier). Equilíbrio ácido-base. Equilíbrios em solução aquosa: solubilidade e íons complexos. " + [char]11 + "Eletroquímica: Semi-reações. Potenciais de eletrodo padrão. Potencial de células galvânicas. Relação entre variação de energia livre padrão e potencial de célula. Energia livre de Gibbs e a Equação de Nernst. Eletrólise e lei de Faraday."

# --- Paragraph 11: Programa resumido body (PT) ---
$d.Paragraphs.Item(11).Range.Text = "Serão realizadas duas provas escritas"

# --- Paragraph 12: Programa resumido body (EN, italic) -> long EN objectives ---
$d.Paragraphs.Item(12).Range.Text = "Make predictions and explain about the effects on rates of reactions in view of: catalysis, temperature variation, collision geometry and concentration of reactants. Suggesting laws of rate using data concerning the effects of concentration. Studying the reactions of order zero, one and two. Considering the application of the Arrhenius equation. Understanding the meaning of spontaneous directions of chemical reactions. Understanding the entropy and its relationship with the three laws of thermodynamics. Understanding the free energy as a reference to the degree of separation between the system and its state of equilibrium. Using the variation of standard free energy as a tool to calculate the equilibrium constant for a given process. Examining the concept of balance and define the equilibrium constant. Learning how to write the expressions of equilibrium constants for homogeneous and heterogeneous reactions. Interpreting the meaning of the reaction quotient. Giving the solubility product expression for a salt, calculate the solubility products from experimentally data of solubility and predict whether precipitation occurs or not. Applying the precipitation criterion. Determining the effects of the common ion. Calculate the pH of an aqueous solution for systems involving strong acid or base or at least a weak acid or a weak base. Understanding the behavior of buffer solutions. Studying the reactions involving the formation of complex ions from high soluble and poor soluble species. Performing a simple representation of a chemical system of an electrolytic or galvanic cell observing the standard conventions to identify the electrodes. Applying the principles of stoichiometry to electrochemical processes using balanced equations of half-reactions and the value of constant of Faraday. Calculate potentials- standard cells using tabulated values. Gibbs free energy and Nernst equation."

# --- Paragraph 14: Programa body (PT detail) -> short criterio formula ---
$d.Paragraphs.Item(14).Range.Text = "NF = (P1 + P2*2)/3"

# --- Paragraph 17: Avaliacao list bullet, replace the 3 non-label runs ---
# Order chosen to avoid transient text collisions between segments.
$rng17 = $d.Paragraphs.Item(17).Range
$rng17.Find.Execute("Será realizada uma avaliação (P3) englobando toda a ementa. A média final será obtida conforme equação: MF= (NF+P3)/2.", $true, $false, $false, $false, $false, $true, 1, $false, "5817330 - Larissa de Freitas", 2) | Out-Null
$rng17 = $d.Paragraphs.Item(17).Range
$rng17.Find.Execute("NF = (P1 + P2*2)/3", $true, $false, $false, $false, $false, $true, 1, $false, "1) ATIKNS, P.; JONES, L. Princípios de Química, 5ªEdição, Ed. Bookman, 2012. 2) BRADY, J.; HUMISTON, G.E. Química Geral Volume II, 2ª Edição, Ed. LTC, 2005. 3) BRADY, J.E.; RUSSELL, J.W.; HOLUM, J.R. Química a matéria e suas transformações Volume II 3ª Edição, Ed. LTC, 2010. 4) BRADY, J.E.; SENESE, F., Química – A matéria e suas transformações Volume II, Ed. LTC, 5ªEdição, 2010. 5) BROWNN, S.L.; HOLME, T.A. Química geral aplicada à engenharia. São Paulo: Ed. Cengage Learning, 2010. 6) BROWN, T. L.; LEMAY, H.E.L.; Jr BURSTEN, B.E.; BURDGE, J.R. Química a ciência central. 9ª Edição, Ed. Pearson Prentice Hall, 2005. 7) CHANG, R.; GOLDSBY, K.A., Química, 11ª Edição, Ed. AMGH Editora Ltda, 2013 8) KOTZ, J.C.; TREICHEL, P.M.; WEAVER, G.C., Química Geral e Reações Químicas, Volume II, 6ª Edição, Ed. Cengage Learning, 2009. 9) KOTZ, J.C.; TREICHEL, P.M.; TOWNSEND, J. R.; TREICHEL, D.A., Química Geral e Reações Químicas, Volume II, 9ª Edição, Ed. Cengage Learning, 2016.", 2) | Out-Null
$rng17 = $d.Paragraphs.Item(17).Range
$rng17.Find.Execute("Serão realizadas duas provas escritas", $true, $false, $false, $false, $false, $true, 1, $false, "Será realizada uma avaliação (P3) englobando toda a ementa. A média final será obtida conforme equação: MF= (NF+P3)/2.", 2) | Out-Null

# --- Paragraph 19: Bibliografia body -> docente name ---
$d.Paragraphs.Item(19).Range.Text = "1506103 - Pedro Carlos de Oliveira"

